# Auto-generated Excel COM-interop edit script
# Applies cached price/profit value updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 45000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45228

$ws.Range("H40").Value = 2333.3333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2333.3333
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2333.3333
$ws.Range("N40").Value = -2683.3333

$ws.Range("H102").Value = 45000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 45000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -51490

$ws.Range("H129").Value = 936.8095
$ws.Range("I129").Value = 349.9
$ws.Range("J129").Value = 1047.5471
$ws.Range("K129").Value = 1049.7
$ws.Range("L129").Value = 3142.6413
$ws.Range("M129").Value = 3950.3
$ws.Range("N129").Value = -13142.6413

$ws.Range("H137").Value = 804.7368
$ws.Range("I137").Value = 749.40625
$ws.Range("J137").Value = 1099.8334
$ws.Range("K137").Value = 2248.21875
$ws.Range("L137").Value = 3299.5002
$ws.Range("M137").Value = 301.78125
$ws.Range("N137").Value = -8399.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1948.3334
$ws.Range("I2").Value = 1743.4615
$ws.Range("J2").Value = 2281.25
$ws.Range("K2").Value = 1743.4615
$ws.Range("L2").Value = 2281.25
$ws.Range("M2").Value = -1630.4615
$ws.Range("N2").Value = -2507.25

$ws.Range("H32").Value = 6653.74
$ws.Range("I32").Value = 6417.109
$ws.Range("J32").Value = 9375
$ws.Range("K32").Value = 6417.109
$ws.Range("L32").Value = 9375
$ws.Range("M32").Value = -6130.109
$ws.Range("N32").Value = -9949

$ws.Range("H61").Value = 1857.409
$ws.Range("I61").Value = 1755.2424
$ws.Range("J61").Value = 2163.9092
$ws.Range("K61").Value = 1755.2424
$ws.Range("L61").Value = 2163.9092
$ws.Range("M61").Value = -1543.2424
$ws.Range("N61").Value = -2587.9092

$ws.Range("H74").Value = 1114.9412
$ws.Range("I74").Value = 1133.6154
$ws.Range("J74").Value = 1054.25
$ws.Range("K74").Value = 1133.6154
$ws.Range("L74").Value = 1054.25
$ws.Range("M74").Value = -259.6153999999999
$ws.Range("N74").Value = -2802.25

$ws.Range("H77").Value = 1114.9412
$ws.Range("I77").Value = 1133.6154
$ws.Range("J77").Value = 1054.25
$ws.Range("K77").Value = 5668.076999999999
$ws.Range("L77").Value = 5271.25
$ws.Range("M77").Value = -1300.076999999999
$ws.Range("N77").Value = -14007.25

$ws.Range("H102").Value = 3339.6155
$ws.Range("I102").Value = 2211.3809
$ws.Range("J102").Value = 8078.2
$ws.Range("K102").Value = 2211.3809
$ws.Range("L102").Value = 8078.2
$ws.Range("M102").Value = -589.3809000000001
$ws.Range("N102").Value = -11322.2

$ws.Range("H109").Value = 52000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 52000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 52000
$ws.Range("N109").Value = -54774

$ws.Range("H116").Value = 1948.3334
$ws.Range("I116").Value = 1743.4615
$ws.Range("J116").Value = 2281.25
$ws.Range("K116").Value = 1743.4615
$ws.Range("L116").Value = 2281.25
$ws.Range("M116").Value = 550.5385000000001
$ws.Range("N116").Value = -6869.25

$ws.Range("H132").Value = 1788.2174
$ws.Range("I132").Value = 1351.625
$ws.Range("J132").Value = 2786.1428
$ws.Range("K132").Value = 4054.875
$ws.Range("L132").Value = 8358.428400000001
$ws.Range("M132").Value = -1524.875
$ws.Range("N132").Value = -13418.4284

$ws.Range("H136").Value = 1857.409
$ws.Range("I136").Value = 1755.2424
$ws.Range("J136").Value = 2163.9092
$ws.Range("K136").Value = 5265.7272
$ws.Range("L136").Value = 6491.7276
$ws.Range("M136").Value = -2715.7272
$ws.Range("N136").Value = -11591.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1948.3334
$ws.Range("I3").Value = 1743.4615
$ws.Range("J3").Value = 2281.25
$ws.Range("K3").Value = 1743.4615
$ws.Range("L3").Value = 2281.25
$ws.Range("M3").Value = -1629.4615
$ws.Range("N3").Value = -2509.25

$ws.Range("H99").Value = 2097.838
$ws.Range("I99").Value = 2742.353
$ws.Range("J99").Value = 1550
$ws.Range("K99").Value = 2742.353
$ws.Range("L99").Value = 1550
$ws.Range("M99").Value = -1244.353
$ws.Range("N99").Value = -4546

$ws.Range("H103").Value = 27000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 27000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 27000
$ws.Range("N103").Value = -29344

$ws.Range("H107").Value = 6541.625
$ws.Range("I107").Value = 1054.5
$ws.Range("J107").Value = 66900
$ws.Range("K107").Value = 1054.5
$ws.Range("L107").Value = 66900
$ws.Range("M107").Value = 865.5
$ws.Range("N107").Value = -70740

$ws.Range("H134").Value = 24751.205
$ws.Range("I134").Value = 1833.9117
$ws.Range("J134").Value = 102670
$ws.Range("K134").Value = 5501.7351
$ws.Range("L134").Value = 308010
$ws.Range("M134").Value = -2966.7351
$ws.Range("N134").Value = -313080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 13326
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 13326
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 13326
$ws.Range("N43").Value = -13694

$ws.Range("H58").Value = 5052.8667
$ws.Range("I58").Value = 1722.091
$ws.Range("J58").Value = 14212.5
$ws.Range("K58").Value = 1722.091
$ws.Range("L58").Value = 14212.5
$ws.Range("M58").Value = -1519.091
$ws.Range("N58").Value = -14618.5

$ws.Range("H75").Value = 38000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 38000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 38000
$ws.Range("N75").Value = -39996

$ws.Range("H78").Value = 38000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 38000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 114000
$ws.Range("N78").Value = -123984

$ws.Range("H101").Value = 13326
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 13326
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 13326
$ws.Range("N101").Value = -19816

$ws.Range("H132").Value = 1741.7778
$ws.Range("I132").Value = 1227.1666
$ws.Range("J132").Value = 2771
$ws.Range("K132").Value = 3681.4998
$ws.Range("L132").Value = 8313
$ws.Range("M132").Value = -1151.4998
$ws.Range("N132").Value = -13373

$ws.Range("H134").Value = 38462664
$ws.Range("I134").Value = 1260
$ws.Range("J134").Value = 166667340
$ws.Range("K134").Value = 3780
$ws.Range("L134").Value = 500002020
$ws.Range("M134").Value = -1245
$ws.Range("N134").Value = -500007090

$ws.Range("H136").Value = 5052.8667
$ws.Range("I136").Value = 1722.091
$ws.Range("J136").Value = 14212.5
$ws.Range("K136").Value = 5166.272999999999
$ws.Range("L136").Value = 42637.5
$ws.Range("M136").Value = -2616.272999999999
$ws.Range("N136").Value = -47737.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1000.11865
$ws.Range("I122").Value = 473.66666
$ws.Range("J122").Value = 1134.5319
$ws.Range("K122").Value = 4262.99994
$ws.Range("L122").Value = 10210.7871
$ws.Range("M122").Value = -1812.99994
$ws.Range("N122").Value = -15110.7871

$ws.Range("H129").Value = 28642.947
$ws.Range("I129").Value = 1013.63635
$ws.Range("J129").Value = 39899.332
$ws.Range("K129").Value = 3040.90905
$ws.Range("L129").Value = 119697.996
$ws.Range("M129").Value = 1959.09095
$ws.Range("N129").Value = -129697.996

$ws.Range("H131").Value = 35859324
$ws.Range("I131").Value = 100202080
$ws.Range("J131").Value = 113349
$ws.Range("K131").Value = 300606240
$ws.Range("L131").Value = 340047
$ws.Range("M131").Value = -300601200
$ws.Range("N131").Value = -350127

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3398.889
$ws.Range("I80").Value = 3450.9524
$ws.Range("J80").Value = 3216.6667
$ws.Range("K80").Value = 3450.9524
$ws.Range("L80").Value = 3216.6667
$ws.Range("M80").Value = -2452.9524
$ws.Range("N80").Value = -5212.6667

$ws.Range("H83").Value = 3398.889
$ws.Range("I83").Value = 3450.9524
$ws.Range("J83").Value = 3216.6667
$ws.Range("K83").Value = 17254.762
$ws.Range("L83").Value = 16083.3335
$ws.Range("M83").Value = -12262.762
$ws.Range("N83").Value = -26067.3335

$ws.Range("H102").Value = 1573.1428
$ws.Range("I102").Value = 1502
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1502
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 120
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2689.1052
$ws.Range("I68").Value = 2741.0833
$ws.Range("J68").Value = 2600
$ws.Range("K68").Value = 2741.0833
$ws.Range("L68").Value = 2600
$ws.Range("M68").Value = -1992.0833
$ws.Range("N68").Value = -4098

$ws.Range("H71").Value = 2689.1052
$ws.Range("I71").Value = 2741.0833
$ws.Range("J71").Value = 2600
$ws.Range("K71").Value = 13705.4165
$ws.Range("L71").Value = 13000
$ws.Range("M71").Value = -9961.416499999999
$ws.Range("N71").Value = -20488

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H93").Value = 2926.8215
$ws.Range("I93").Value = 3219.0557
$ws.Range("J93").Value = 2400.8
$ws.Range("K93").Value = 3219.0557
$ws.Range("L93").Value = 2400.8
$ws.Range("M93").Value = -1971.0557
$ws.Range("N93").Value = -4896.8

$ws.Range("H136").Value = 3608.9348
$ws.Range("I136").Value = 1894.8948
$ws.Range("J136").Value = 11750.625
$ws.Range("K136").Value = 5684.6844
$ws.Range("L136").Value = 35251.875
$ws.Range("M136").Value = -3134.6844
$ws.Range("N136").Value = -40351.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 24500
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 24500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 24500
$ws.Range("N27").Value = -24638

$ws.Range("H107").Value = 509.9697
$ws.Range("I107").Value = 445.72415
$ws.Range("J107").Value = 975.75
$ws.Range("K107").Value = 1337.17245
$ws.Range("L107").Value = 2927.25
$ws.Range("M107").Value = 582.82755
$ws.Range("N107").Value = -6767.25

$ws.Range("H109").Value = 11562.375
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 11562.375
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 11562.375
$ws.Range("N109").Value = -14336.375

$ws.Range("H132").Value = 1081.5098
$ws.Range("I132").Value = 860.3095
$ws.Range("J132").Value = 2113.7778
$ws.Range("K132").Value = 2580.9285
$ws.Range("L132").Value = 6341.3334
$ws.Range("M132").Value = -50.92849999999999
$ws.Range("N132").Value = -11401.3334

Write-Host "Applied $([int]44) row updates across 8 sheets"
